$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("June 2018")

$xlCenter = -4108

# Fix completion date for "Chest Reopening Concept 3 Left + Right Sternum" (row 30, column B = Date Completed)
$ws.Cells.Item(30, 2).Value = "29-06-2018"

# Add new request row 31: 200 Tracheal Rings
$ws.Cells.Item(31, 1).Value = "29-06-2018"
$ws.Cells.Item(31, 3).Value = "Tracheal Rings od: 23.5 id: 15.5"
$ws.Cells.Item(31, 4).Value = 200
$ws.Cells.Item(31, 5).Value = "Polylite/PLA"
$ws.Cells.Item(31, 6).Value = 2
$ws.Cells.Item(31, 7).Value = 20
$ws.Cells.Item(31, 8).Value = 0.2
$ws.Cells.Item(31, 9).Value = "NA"

# Match formatting (centered alignment, same style as rest of the data rows) used by other rows
$ws.Cells.Item(31, 1).HorizontalAlignment = $xlCenter
$ws.Cells.Item(31, 3).HorizontalAlignment = $xlCenter
$ws.Cells.Item(31, 4).HorizontalAlignment = $xlCenter
$ws.Cells.Item(31, 5).HorizontalAlignment = $xlCenter
$ws.Cells.Item(31, 6).HorizontalAlignment = $xlCenter
$ws.Cells.Item(31, 7).HorizontalAlignment = $xlCenter
$ws.Cells.Item(31, 8).HorizontalAlignment = $xlCenter
$ws.Cells.Item(31, 9).HorizontalAlignment = $xlCenter

# Update selection to reflect the new active cell on this sheet
$ws.Range("A31").Select()
